$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the last data row (row 6); dataset trimmed by one reading ---
$ws.Rows.Item(6).Delete()

# --- 2. Widen several data columns by one character unit (7 -> 8, except col 20: 8 -> 9) ---
# COM ColumnWidth uses the workbook default font metrics; this engine maps
# OOXML <col width="N"> = ColumnWidth + 0.8333333333333334, so we back that offset out
# to land exactly on the integer widths recorded in the target file.
$colWidthOffset = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 8 - $colWidthOffset   # B -> 8
$ws.Columns.Item(3).ColumnWidth = 8 - $colWidthOffset   # C -> 8
$ws.Columns.Item(7).ColumnWidth = 8 - $colWidthOffset   # G -> 8
$ws.Columns.Item(9).ColumnWidth = 8 - $colWidthOffset   # I -> 8
$ws.Columns.Item(10).ColumnWidth = 8 - $colWidthOffset   # J -> 8
$ws.Columns.Item(11).ColumnWidth = 8 - $colWidthOffset   # K -> 8
$ws.Columns.Item(12).ColumnWidth = 8 - $colWidthOffset   # L -> 8
$ws.Columns.Item(13).ColumnWidth = 8 - $colWidthOffset   # M -> 8
$ws.Columns.Item(15).ColumnWidth = 8 - $colWidthOffset   # O -> 8
$ws.Columns.Item(16).ColumnWidth = 8 - $colWidthOffset   # P -> 8
$ws.Columns.Item(17).ColumnWidth = 8 - $colWidthOffset   # Q -> 8
$ws.Columns.Item(20).ColumnWidth = 9 - $colWidthOffset   # T -> 9
$ws.Columns.Item(22).ColumnWidth = 8 - $colWidthOffset   # V -> 8
$ws.Columns.Item(23).ColumnWidth = 8 - $colWidthOffset   # W -> 8
$ws.Columns.Item(24).ColumnWidth = 8 - $colWidthOffset   # X -> 8
$ws.Columns.Item(26).ColumnWidth = 8 - $colWidthOffset   # Z -> 8
$ws.Columns.Item(27).ColumnWidth = 8 - $colWidthOffset   # AA -> 8
$ws.Columns.Item(28).ColumnWidth = 8 - $colWidthOffset   # AB -> 8
$ws.Columns.Item(29).ColumnWidth = 8 - $colWidthOffset   # AC -> 8
$ws.Columns.Item(30).ColumnWidth = 8 - $colWidthOffset   # AD -> 8
$ws.Columns.Item(34).ColumnWidth = 8 - $colWidthOffset   # AH -> 8

# --- 3. Replace the 4 remaining data rows (2-5) with the new reading values ---

# Row 2
$ws.Cells.Item(2, 1).Value = 45068.50694444445
$ws.Cells.Item(2, 2).Value = 21.619
$ws.Cells.Item(2, 3).Value = 14.458
$ws.Cells.Item(2, 4).Value = 4.456
$ws.Cells.Item(2, 5).Value = 45.361
$ws.Cells.Item(2, 6).Value = 37.513
$ws.Cells.Item(2, 7).Value = 17.013
$ws.Cells.Item(2, 8).Value = 55.703
$ws.Cells.Item(2, 9).Value = 26.178
$ws.Cells.Item(2, 10).Value = 11.065
$ws.Cells.Item(2, 11).Value = 17.024
$ws.Cells.Item(2, 12).Value = 18.048
$ws.Cells.Item(2, 13).Value = 18.877
$ws.Cells.Item(2, 14).Value = 5.432
$ws.Cells.Item(2, 15).Value = 16.918
$ws.Cells.Item(2, 16).Value = 23.745
$ws.Cells.Item(2, 17).Value = 14.242
$ws.Cells.Item(2, 18).Value = 3.744
$ws.Cells.Item(2, 19).Value = 2.457
$ws.Cells.Item(2, 20).Value = 250.169
$ws.Cells.Item(2, 21).Value = 47.091
$ws.Cells.Item(2, 22).Value = 15.616
$ws.Cells.Item(2, 23).Value = 31.148
$ws.Cells.Item(2, 24).Value = 16.247
$ws.Cells.Item(2, 25).Value = 2.407
$ws.Cells.Item(2, 26).Value = 27.662
$ws.Cells.Item(2, 27).Value = 13.794
$ws.Cells.Item(2, 28).Value = 12.753
$ws.Cells.Item(2, 29).Value = 14.55
$ws.Cells.Item(2, 30).Value = 18.409
$ws.Cells.Item(2, 31).Value = 3.64
$ws.Cells.Item(2, 32).Value = 49.22
$ws.Cells.Item(2, 33).Value = 8.641
$ws.Cells.Item(2, 34).Value = 19.523

# Row 3
$ws.Cells.Item(3, 1).Value = 45068.51388888889
$ws.Cells.Item(3, 2).Value = 10.089
$ws.Cells.Item(3, 3).Value = 6.826
$ws.Cells.Item(3, 4).Value = 1.64
$ws.Cells.Item(3, 5).Value = 21.283
$ws.Cells.Item(3, 6).Value = 17.604
$ws.Cells.Item(3, 7).Value = 7.94
$ws.Cells.Item(3, 8).Value = 33.495
$ws.Cells.Item(3, 9).Value = 12.216
$ws.Cells.Item(3, 10).Value = 5.176
$ws.Cells.Item(3, 11).Value = 7.827
$ws.Cells.Item(3, 12).Value = 8.614000000000001
$ws.Cells.Item(3, 13).Value = 8.913
$ws.Cells.Item(3, 14).Value = 2.539
$ws.Cells.Item(3, 15).Value = 7.895
$ws.Cells.Item(3, 16).Value = 11.072
$ws.Cells.Item(3, 17).Value = 6.89
$ws.Cells.Item(3, 18).Value = 1.56
$ws.Cells.Item(3, 19).Value = 0.884
$ws.Cells.Item(3, 20).Value = 112.867
$ws.Cells.Item(3, 21).Value = 22.214
$ws.Cells.Item(3, 22).Value = 7.288
$ws.Cells.Item(3, 23).Value = 14.571
$ws.Cells.Item(3, 24).Value = 7.818
$ws.Cells.Item(3, 25).Value = 1.182
$ws.Cells.Item(3, 26).Value = 15.689
$ws.Cells.Item(3, 27).Value = 6.437
$ws.Cells.Item(3, 28).Value = 5.997
$ws.Cells.Item(3, 29).Value = 6.911
$ws.Cells.Item(3, 30).Value = 8.865
$ws.Cells.Item(3, 31).Value = 1.294
$ws.Cells.Item(3, 32).Value = 30.477
$ws.Cells.Item(3, 33).Value = 3.967
$ws.Cells.Item(3, 34).Value = 9.112

# Row 4
$ws.Cells.Item(4, 1).Value = 45068.52083333334
$ws.Cells.Item(4, 2).Value = 21.619
$ws.Cells.Item(4, 3).Value = 15.762
$ws.Cells.Item(4, 4).Value = 1.511
$ws.Cells.Item(4, 5).Value = 46.599
$ws.Cells.Item(4, 6).Value = 38.53
$ws.Cells.Item(4, 7).Value = 17.013
$ws.Cells.Item(4, 8).Value = 63.837
$ws.Cells.Item(4, 9).Value = 26.178
$ws.Cells.Item(4, 10).Value = 11.543
$ws.Cells.Item(4, 11).Value = 17.281
$ws.Cells.Item(4, 12).Value = 18.804
$ws.Cells.Item(4, 13).Value = 19.722
$ws.Cells.Item(4, 14).Value = 5.434
$ws.Cells.Item(4, 15).Value = 16.918
$ws.Cells.Item(4, 16).Value = 24.015
$ws.Cells.Item(4, 17).Value = 14.307
$ws.Cells.Item(4, 18).Value = 1.122
$ws.Cells.Item(4, 19).Value = 0.97
$ws.Cells.Item(4, 20).Value = 250.194
$ws.Cells.Item(4, 21).Value = 47.262
$ws.Cells.Item(4, 22).Value = 15.616
$ws.Cells.Item(4, 23).Value = 31.672
$ws.Cells.Item(4, 24).Value = 16.845
$ws.Cells.Item(4, 25).Value = 2.338
$ws.Cells.Item(4, 26).Value = 31.14
$ws.Cells.Item(4, 27).Value = 13.794
$ws.Cells.Item(4, 28).Value = 12.329
$ws.Cells.Item(4, 29).Value = 14.436
$ws.Cells.Item(4, 30).Value = 19.654
$ws.Cells.Item(4, 31).Value = 0.784
$ws.Cells.Item(4, 32).Value = 57.73
$ws.Cells.Item(4, 33).Value = 8.759
$ws.Cells.Item(4, 34).Value = 19.524

# Row 5
$ws.Cells.Item(5, 1).Value = 45068.52777777778
$ws.Cells.Item(5, 2).Value = 10.57
$ws.Cells.Item(5, 3).Value = 7.59
$ws.Cells.Item(5, 4).Value = 0.92
$ws.Cells.Item(5, 5).Value = 22.67
$ws.Cells.Item(5, 6).Value = 18.74
$ws.Cells.Item(5, 7).Value = 8.32
$ws.Cells.Item(5, 8).Value = 35.67
$ws.Cells.Item(5, 9).Value = 12.8
$ws.Cells.Item(5, 10).Value = 5.61
$ws.Cells.Item(5, 11).Value = 8.359999999999999
$ws.Cells.Item(5, 12).Value = 9.199999999999999
$ws.Cells.Item(5, 13).Value = 9.59
$ws.Cells.Item(5, 14).Value = 2.66
$ws.Cells.Item(5, 15).Value = 8.27
$ws.Cells.Item(5, 16).Value = 11.72
$ws.Cells.Item(5, 17).Value = 7.1
$ws.Cells.Item(5, 18).Value = 0.78
$ws.Cells.Item(5, 19).Value = 0.54
$ws.Cells.Item(5, 20).Value = 118.59
$ws.Cells.Item(5, 21).Value = 23.26
$ws.Cells.Item(5, 22).Value = 7.63
$ws.Cells.Item(5, 23).Value = 15.48
$ws.Cells.Item(5, 24).Value = 8.279999999999999
$ws.Cells.Item(5, 25).Value = 1.19
$ws.Cells.Item(5, 26).Value = 16.73
$ws.Cells.Item(5, 27).Value = 6.74
$ws.Cells.Item(5, 28).Value = 6.09
$ws.Cells.Item(5, 29).Value = 7.11
$ws.Cells.Item(5, 30).Value = 9.59
$ws.Cells.Item(5, 31).Value = 0.5600000000000001
$ws.Cells.Item(5, 32).Value = 32.5
$ws.Cells.Item(5, 33).Value = 4.23
$ws.Cells.Item(5, 34).Value = 9.550000000000001
